# Generate Report for Handoff
#
# "b.md" has finished its handoff cycle (source
# b.63290e5768f688058c7b37413b0a5c26c308f864 was handed off for both the
# zh-cn and de-de locales), so its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", the "Latest
# Handoff File" / "Latest Handoff Datetime" columns are refreshed to point
# at the new .xlf, and the Overview sheet's handoff date is refreshed to
# match.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = "2016-24-17 12:24:40"

# ---------------------------------------------------------------------------
# zh-cn sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-17 12:24:37"

# ---------------------------------------------------------------------------
# de-de sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-17 12:24:40"

$wb.Save()
